$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Header row: lower-case / snake_case the column headers, and add
#    two new trailing headers (K1, L1) for the columns that the
#    Probable/Hospital block needs once it is shifted right.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "confirmed"
$ws.Range("C1").Value = "total_confirmed"
$ws.Range("D1").Value = "probable"
$ws.Range("E1").Value = "total_probable"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "cumulative_total"
$ws.Range("H1").Value = "recovered"
$ws.Range("I1").Value = "total_recovered"
$ws.Range("J1").Value = "in_hospital_now"
$ws.Range("K1").Value = "total_been_in_hospital"
$ws.Range("L1").Value = "in_icu"

# ------------------------------------------------------------------
# 2. Rows 2-26: previously only had A (date), B (confirmed) and
#    C (total_confirmed). Add the new "total" (F) and
#    "cumulative_total" (G) daily recovered-style running counts.
# ------------------------------------------------------------------
$fg = @{
  2  = @(1, 1)
  3  = @(1, 2)
  4  = @(1, 3)
  5  = @(1, 4)
  6  = @(1, 5)
  7  = @(1, 6)
  8  = @(1, 7)
  9  = @(1, 8)
  10 = @(1, 9)
  11 = @(1, 10)
  12 = @(1, 11)
  13 = @(1, 12)
  14 = @(1, 13)
  15 = @(1, 14)
  16 = @(1, 15)
  17 = @(1, 16)
  18 = @(2, 18)
  19 = @(2, 20)
  20 = @(4, 24)
  21 = @(8, 32)
  22 = @(8, 40)
  23 = @(11, 51)
  24 = @(13, 64)
  25 = @(14, 78)
  26 = @(36, 114)
}

foreach ($r in $fg.Keys) {
  $vals = $fg[$r]
  $ws.Cells.Item($r, 6).Value = $vals[0]   # F = total
  $ws.Cells.Item($r, 7).Value = $vals[1]   # G = cumulative_total
}

# ------------------------------------------------------------------
# 3. Rows 27-30: these already had data out through column J
#    (Probable/TotalProbable/Recovered/TotalRecovered/InHospitalNow/
#    TotalBeenInHospital/InICU in the old layout). The new MoH layout
#    inserts "total" + "cumulative_total" after total_probable (E),
#    shifting the old F..J block to H..L. Rewrite each row completely
#    in the new column order.
# ------------------------------------------------------------------
function Set-Row($r, $vals) {
  # $vals is an ordered hashtable-like array of [col, value] pairs
  foreach ($pair in $vals) {
    $ws.Cells.Item($r, $pair[0]).Value = $pair[1]
  }
}

# row 27
Set-Row 27 @(
  ,(1, 43914)
  ,(2, 40)
  ,(3, 142)
  ,(4, 13)
  ,(5, 13)
  ,(6, 53)
  ,(7, 167)
  ,(8, 12)
  ,(9, 12)
  ,(10, 6)
  ,(12, 0)
)

# row 28
Set-Row 28 @(
  ,(1, 43915)
  ,(2, 47)
  ,(3, 189)
  ,(4, 3)
  ,(5, 16)
  ,(6, 50)
  ,(7, 217)
  ,(8, 10)
  ,(9, 22)
  ,(10, 6)
  ,(12, 0)
)

# row 29
Set-Row 29 @(
  ,(1, 43916)
  ,(2, 73)
  ,(3, 262)
  ,(4, 5)
  ,(5, 21)
  ,(6, 78)
  ,(7, 295)
  ,(8, 5)
  ,(9, 27)
  ,(10, 7)
  ,(12, 0)
)

# row 30
Set-Row 30 @(
  ,(1, 43917)
  ,(2, 76)
  ,(3, 338)
  ,(4, 9)
  ,(5, 30)
  ,(6, 85)
  ,(7, 380)
  ,(8, 10)
  ,(9, 37)
  ,(10, 8)
  ,(11, 20)
  ,(12, 1)
)

# ------------------------------------------------------------------
# 4. New row 31 with the latest day's figures, same full layout.
# ------------------------------------------------------------------
Set-Row 31 @(
  ,(1, 43918)
  ,(2, 78)
  ,(3, 416)
  ,(4, 5)
  ,(5, 35)
  ,(6, 83)
  ,(7, 463)
  ,(8, 13)
  ,(9, 50)
  ,(10, 12)
  ,(11, 22)
  ,(12, 2)
)

# Row 31's date cell should carry the same date-number-format style (s="2")
# that all the other date cells in column A use.
$ws.Cells.Item(31, 1).NumberFormat = "yyyy-mm-dd"
